$wb = $excel.ActiveWorkbook

# A new handback cycle completed for the
# 9eb30ea2-3d8e-422e-a4af-de19f77b8121 file (row 7 of both language
# sheets). Refresh its Correspond Handoff Datetime (D) and Correspond
# Handback DateTime (G) columns for each locale's report sheet.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D7").Value = "2016-03-03 10:25:26"
$wsZhCn.Range("G7").Value = "2016-03-03 10:26:09"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D7").Value = "2016-03-03 10:25:37"
$wsDeDe.Range("G7").Value = "2016-03-03 10:26:32"
